$d = $word.ActiveDocument

# Locate "-2020" (the hyphen immediately followed by "2020") inside the
# "Data Analytics Bootcamp, 2019-2020" heading line.
$findRange = $d.Content
$found = $findRange.Find.Execute("-2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the '-2020' text to update."
}

# Re-anchor to a fresh Range over the same span so InsertXML replaces
# exactly (and only) the matched text.
$target = $d.Range($findRange.Start, $findRange.End)

# Replace the found span with " - 2020", split across four runs (matching
# how Word itself records such an incremental retype) while keeping the
# same run formatting (bold, majorHAnsi theme fonts, noProof) used by the
# surrounding heading text.
$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:noProof/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:noProof/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:noProof/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:noProof/></w:rPr><w:t>2020</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$target.InsertXML($xml)
